$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "'-0.92%"
$ws.Range("D3").Value = "'31.08"
$ws.Range("E3").Value = "'1.33%"
$ws.Range("D4").Value = "'4.930"
$ws.Range("E4").Value = "'-0.51%"
$ws.Range("D5").Value = "'0.07339"
$ws.Range("E5").Value = "'1.81%"
$ws.Range("D6").Value = "'2.230"
$ws.Range("E6").Value = "'23.53%"
$ws.Range("D7").Value = "'7.718"
$ws.Range("E7").Value = "'0.52%"
$ws.Range("D8").Value = "'3.726"
$ws.Range("E8").Value = "'-1.12%"
$ws.Range("D9").Value = "'0.9017"
$ws.Range("E9").Value = "'0.57%"
$ws.Range("D10").Value = "'0.09198"
$ws.Range("E10").Value = "'18.67%"
$ws.Range("D11").Value = "'0.1694"
$ws.Range("E11").Value = "'2.31%"
$ws.Range("D12").Value = "'0.08217"
$ws.Range("E12").Value = "'2.16%"
$ws.Range("E13").Value = "'2.31%"
$ws.Range("D14").Value = "'0.09940"
$ws.Range("E14").Value = "'-0.67%"
$ws.Range("D15").Value = "'0.001496"
$ws.Range("E15").Value = "'-0.11%"
$ws.Range("D16").Value = "'0.005732"
$ws.Range("E16").Value = "'0.75%"
$ws.Range("D17").Value = "'3.525"
$ws.Range("E17").Value = "'1.71%"
$ws.Range("D18").Value = "'2.073"
$ws.Range("E18").Value = "'-0.34%"
$ws.Range("D19").Value = "'0.3332"
$ws.Range("E19").Value = "'0.47%"
$ws.Range("E20").Value = "'0.02%"
$ws.Range("D21").Value = "'4.155"
$ws.Range("E21").Value = "'2.81%"
$ws.Range("E22").Value = "'-12.07%"
$ws.Range("D23").Value = "'0.04532"
$ws.Range("E23").Value = "'0.50%"
$ws.Range("D24").Value = "'0.001209"
$ws.Range("E24").Value = "'-0.61%"
$ws.Range("D25").Value = "'0.004165"
$ws.Range("E25").Value = "'4.15%"
$ws.Range("E26").Value = "'3.89%"
$ws.Range("D27").Value = "'0.0003395"
$ws.Range("D39").Value = "'0.01568"
$ws.Range("E39").Value = "'-0.76%"
$ws.Range("D40").Value = "'0.04442"
$ws.Range("E40").Value = "'0.80%"
$ws.Range("D41").Value = "'0.007341"
$ws.Range("E41").Value = "'1.28%"
$ws.Range("D42").Value = "'0.009543"
$ws.Range("E42").Value = "'-3.54%"
$ws.Range("E43").Value = "'1.85%"
$ws.Range("D44").Value = "'0.002310"
$ws.Range("E44").Value = "'15.02%"
$ws.Range("D45").Value = "'0.009078"
$ws.Range("E45").Value = "'-4.51%"
$ws.Range("D46").Value = "'0.00006116"
$ws.Range("E46").Value = "'1.98%"
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("E47").Value = "'-0.12%"
$ws.Range("D48").Value = "'2.364"
$ws.Range("E48").Value = "'5.25%"
$ws.Range("E49").Value = "'-33.36%"
$ws.Range("D50").Value = "'0.00002100"
$ws.Range("E50").Value = "'-0.12%"
$ws.Range("D51").Value = "'0.0002000"
$ws.Range("E51").Value = "'-0.12%"
